# Management_information_DKI_2019.xlsx update:
# "Additional translation, standardization and missing data to management files."
#
# Content changes (rows 47-51 of the Managements table):
#  - E47: fertilizer agent "KAS" -> standardized "CAN"
#  - E48: fertilizer agent "Yara Vera Ami Plus" -> standardized "UAN (inhibited)"
#  - B49: stray leftover label "3rd nitrogen application" cleared (no 3rd application recorded)
#  - B50/C50/D50/E50: stray leftover "others" / date / amount / "Bittersalz" cleared
#  - C51/D51/E51: stray leftover date / amount / "Bittersalz" cleared
#
# These were residual placeholder values that no longer apply, so the cells
# are blanked while keeping their existing number formatting/styles intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standardize fertilizer agent names
$ws.Range("E47").Value = "CAN"
$ws.Range("E48").Value = "UAN (inhibited)"

# Clear stale / incorrect leftover entries
$ws.Range("B49").Value = ""

$ws.Range("B50").Value = ""
$ws.Range("C50").Value = ""
$ws.Range("D50").Value = ""
$ws.Range("E50").Value = ""

$ws.Range("C51").Value = ""
$ws.Range("D51").Value = ""
$ws.Range("E51").Value = ""

# Match the author's final selection/view state on save
$ws.Range("E49").Select()
